$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference "plain" style (no font/fill override) taken from A1,
# used to strip the old font-2 / fill-2 based highlighting from B1, C1, J1, A4, A10.
$plainStyle = $ws.Range("A1").Style

# Row 1
$ws.Range("A1").Style = $plainStyle
$ws.Range("B1").Value = 6
$ws.Range("B1").Style = $plainStyle
$ws.Range("C1").Value = 3
$ws.Range("C1").Style = $plainStyle
$ws.Range("D1").Value = 1.5
$ws.Range("D1").Style = $plainStyle
$ws.Range("E1").Value = 0.75
$ws.Range("E1").Style = $plainStyle
$ws.Range("F1").Value = 0.375
$ws.Range("F1").Style = $plainStyle
$ws.Range("G1").Value = 0.1875
$ws.Range("G1").Style = $plainStyle
$ws.Range("H1").Value = 0.09375
$ws.Range("H1").Style = $plainStyle
$ws.Range("I1").Value = 0.046875
$ws.Range("I1").Style = $plainStyle
$ws.Range("J1").Value = 0
$ws.Range("J1").Style = $plainStyle

# Row 2
$ws.Range("A2").Value = 20
$ws.Range("A2").Style = $plainStyle
$ws.Range("B2").Value = 0.99920127795527147
$ws.Range("B2").Style = $plainStyle
$ws.Range("C2").Value = 1.0007987220447283
$ws.Range("C2").Style = $plainStyle
$ws.Range("D2").Value = 0.99920127795527147
$ws.Range("D2").Style = $plainStyle
$ws.Range("E2").Value = 0.99840255591054317
$ws.Range("E2").Style = $plainStyle
$ws.Range("F2").Value = 1
$ws.Range("F2").Style = $plainStyle
$ws.Range("G2").Value = 0.99920127795527147
$ws.Range("G2").Style = $plainStyle
$ws.Range("H2").Value = 1
$ws.Range("H2").Style = $plainStyle
$ws.Range("I2").Value = 1.0007987220447283
$ws.Range("I2").Style = $plainStyle
$ws.Range("J2").Value = 0.99760383386581464
$ws.Range("J2").Style = $plainStyle

# Row 3
$ws.Range("A3").Value = 10
$ws.Range("A3").Style = $plainStyle
$ws.Range("B3").Value = 1.0007987220447283
$ws.Range("B3").Style = $plainStyle
$ws.Range("C3").Value = 1.0023961661341851
$ws.Range("C3").Style = $plainStyle
$ws.Range("D3").Value = 1.0007987220447283
$ws.Range("D3").Style = $plainStyle
$ws.Range("E3").Value = 0.99920127795527147
$ws.Range("E3").Style = $plainStyle
$ws.Range("F3").Value = 1.0007987220447283
$ws.Range("F3").Style = $plainStyle
$ws.Range("G3").Value = 1.0015974440894568
$ws.Range("G3").Style = $plainStyle
$ws.Range("H3").Value = 1
$ws.Range("H3").Style = $plainStyle
$ws.Range("I3").Value = 1.0007987220447283
$ws.Range("I3").Style = $plainStyle
$ws.Range("J3").Value = 0.97364217252396157
$ws.Range("J3").Style = $plainStyle

# Row 4
$ws.Range("A4").Value = 5
$ws.Range("A4").Style = $plainStyle
$ws.Range("B4").Value = 1
$ws.Range("B4").Style = $plainStyle
$ws.Range("C4").Value = 1.0007987220447283
$ws.Range("C4").Style = $plainStyle
$ws.Range("D4").Value = 1
$ws.Range("D4").Style = $plainStyle
$ws.Range("E4").Value = 0.99520766773162939
$ws.Range("E4").Style = $plainStyle
$ws.Range("F4").Value = 1.0015974440894568
$ws.Range("F4").Style = $plainStyle
$ws.Range("G4").Value = 0.99920127795527147
$ws.Range("G4").Style = $plainStyle
$ws.Range("H4").Value = 0.99
$ws.Range("H4").Style = $plainStyle
$ws.Range("I4").Value = 0.99520766773162939
$ws.Range("I4").Style = $plainStyle
$ws.Range("J4").Value = 0.59424920127795522
$ws.Range("J4").Style = $plainStyle

# Row 5
$ws.Range("A5").Value = 2.5
$ws.Range("A5").Style = $plainStyle
$ws.Range("B5").Value = 1.0007987220447283
$ws.Range("B5").Style = $plainStyle
$ws.Range("C5").Value = 1.0015974440894568
$ws.Range("C5").Style = $plainStyle
$ws.Range("D5").Value = 1
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = 1
$ws.Range("E5").Style = $plainStyle
$ws.Range("F5").Value = 1.0007987220447283
$ws.Range("F5").Style = $plainStyle
$ws.Range("G5").Value = 1.0007987220447283
$ws.Range("G5").Style = $plainStyle
$ws.Range("H5").Value = 0.50239616613418525
$ws.Range("H5").Style = $plainStyle
$ws.Range("I5").Value = 0.43370607028753988
$ws.Range("I5").Style = $plainStyle
$ws.Range("J5").Value = -0.069488817891373775
$ws.Range("J5").Style = $plainStyle

# Row 6
$ws.Range("A6").Value = 1.25
$ws.Range("A6").Style = $plainStyle
$ws.Range("B6").Value = 0.99920127795527147
$ws.Range("B6").Style = $plainStyle
$ws.Range("C6").Value = 1.0015974440894568
$ws.Range("C6").Style = $plainStyle
$ws.Range("D6").Value = 0.99201277955271561
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = 0.99920127795527147
$ws.Range("E6").Style = $plainStyle
$ws.Range("F6").Value = 0.99920127795527147
$ws.Range("F6").Style = $plainStyle
$ws.Range("G6").Value = 0.81549520766773154
$ws.Range("G6").Style = $plainStyle
$ws.Range("H6").Value = 0.58706070287539924
$ws.Range("H6").Style = $plainStyle
$ws.Range("I6").Value = 0.24440894568690091
$ws.Range("I6").Style = $plainStyle
$ws.Range("J6").Value = -0.36980830670926523
$ws.Range("J6").Style = $plainStyle

# Row 7
$ws.Range("A7").Value = 0.625
$ws.Range("A7").Style = $plainStyle
$ws.Range("B7").Value = 1.0007987220447283
$ws.Range("B7").Style = $plainStyle
$ws.Range("C7").Value = 1.0023961661341851
$ws.Range("C7").Style = $plainStyle
$ws.Range("D7").Value = 1.0007987220447283
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = 0.58706070287539924
$ws.Range("E7").Style = $plainStyle
$ws.Range("F7").Value = 0.064696485623003161
$ws.Range("F7").Style = $plainStyle
$ws.Range("G7").Value = 0.012779552715654964
$ws.Range("G7").Style = $plainStyle
$ws.Range("H7").Value = -0.22444089456869021
$ws.Range("H7").Style = $plainStyle
$ws.Range("I7").Value = -0.32268370607028762
$ws.Range("I7").Style = $plainStyle
$ws.Range("J7").Value = -0.52316293929712465
$ws.Range("J7").Style = $plainStyle

# Row 8
$ws.Range("A8").Value = 0.3125
$ws.Range("A8").Style = $plainStyle
$ws.Range("B8").Value = 1
$ws.Range("B8").Style = $plainStyle
$ws.Range("C8").Value = 1.0015974440894568
$ws.Range("C8").Style = $plainStyle
$ws.Range("D8").Value = 0.99920127795527147
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = -0.083865814696485602
$ws.Range("E8").Style = $plainStyle
$ws.Range("F8").Value = 0.014376996805111834
$ws.Range("F8").Style = $plainStyle
$ws.Range("G8").Value = -0.13258785942492024
$ws.Range("G8").Style = $plainStyle
$ws.Range("H8").Value = -0.60383386581469645
$ws.Range("H8").Style = $plainStyle
$ws.Range("I8").Value = -0.58067092651757191
$ws.Range("I8").Style = $plainStyle
$ws.Range("J8").Value = -0.74361022364217277
$ws.Range("J8").Style = $plainStyle

# Row 9
$ws.Range("A9").Value = 0.15625
$ws.Range("A9").Style = $plainStyle
$ws.Range("B9").Value = 1.0015974440894568
$ws.Range("B9").Style = $plainStyle
$ws.Range("C9").Value = 1.0023961661341851
$ws.Range("C9").Style = $plainStyle
$ws.Range("D9").Value = 1
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = 0.021565495207667661
$ws.Range("E9").Style = $plainStyle
$ws.Range("F9").Value = -0.44089456869009586
$ws.Range("F9").Style = $plainStyle
$ws.Range("G9").Value = -0.36022364217252401
$ws.Range("G9").Style = $plainStyle
$ws.Range("H9").Value = -0.49920127795527158
$ws.Range("H9").Style = $plainStyle
$ws.Range("I9").Value = -0.47763578274760388
$ws.Range("I9").Style = $plainStyle
$ws.Range("J9").Value = -0.66134185303514381
$ws.Range("J9").Style = $plainStyle

# Row 10
$ws.Range("A10").Value = 0
$ws.Range("A10").Style = $plainStyle
$ws.Range("B10").Value = 1.0007987220447283
$ws.Range("B10").Style = $plainStyle
$ws.Range("C10").Value = 1
$ws.Range("C10").Style = $plainStyle
$ws.Range("D10").Value = 1
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = 1
$ws.Range("E10").Style = $plainStyle
$ws.Range("F10").Value = 0.96565495207667718
$ws.Range("F10").Style = $plainStyle
$ws.Range("G10").Value = 0.0071884984025558278
$ws.Range("G10").Style = $plainStyle
$ws.Range("H10").Value = -0.30670926517571895
$ws.Range("H10").Style = $plainStyle
$ws.Range("I10").Value = -0.34025559105431324
$ws.Range("I10").Style = $plainStyle
$ws.Range("J10").Value = 0
$ws.Range("J10").Style = $plainStyle

# Move the active selection to H18 (matches the saved sheet view state).
$ws.Range("H18").Select() | Out-Null
